$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.674.96"
$ws.Range("E2").Value = "  +1.72%  "

$ws.Range("D3").Value = "1.630.08"
$ws.Range("E3").Value = "  +1.85%  "

$ws.Range("E4").Value = "  -0.30%  "

$ws.Range("D5").Value = "'213.37"
$ws.Range("E5").Value = "  +0.59%  "

$ws.Range("E6").Value = "  -0.28%  "

$ws.Range("E7").Value = "  +0.89%  "

$ws.Range("E8").Value = "  +0.79%  "

$ws.Range("D9").Value = "'0.0620"
$ws.Range("E9").Value = "  +1.15%  "

$ws.Range("D10").Value = "'19.10"
$ws.Range("E10").Value = "  +5.68%  "

$ws.Range("D11").Value = "'0.0828"
$ws.Range("E11").Value = "  +1.96%  "

$ws.Range("D12").Value = "1.855.91"
$ws.Range("E12").Value = "  +1.80%  "

$ws.Range("D13").Value = "1.604.04"
$ws.Range("E13").Value = "  +0.20%  "

$ws.Range("D14").Value = "'4.06"
$ws.Range("E14").Value = "  +0.57%  "

$ws.Range("D15").Value = "'0.524"
$ws.Range("E15").Value = "  +2.29%  "

$ws.Range("D16").Value = "26.625.69"
$ws.Range("E16").Value = "  +1.60%  "

$ws.Range("D17").Value = "'63.10"
$ws.Range("E17").Value = "  +3.09%  "

$ws.Range("D18").Value = "0.0₃0734"
$ws.Range("E18").Value = "  +0.58%  "

$ws.Range("E19").Value = "  -0.22%  "

$ws.Range("D20").Value = "'206.60"
$ws.Range("E20").Value = "  +2.33%  "

$ws.Range("E21").Value = "  +1.16%  "

$ws.Range("D22").Value = "'9.42"
$ws.Range("E22").Value = "  +1.65%  "

$ws.Range("D23").Value = "'6.09"
$ws.Range("E23").Value = "  +1.66%  "

$ws.Range("D24").Value = "'1.89"
$ws.Range("E24").Value = "  -1.94%  "

$ws.Range("D25").Value = "'145.32"
$ws.Range("E25").Value = "  +0.51%  "

$ws.Range("E26").Value = "  -0.33%  "

$ws.Range("E27").Value = "  -0.74%  "

$ws.Range("D28").Value = "'15.43"
$ws.Range("E28").Value = "  +1.76%  "

$ws.Range("E29").Value = "  +1.59%  "

$ws.Range("D30").Value = "'0.0522"
$ws.Range("E30").Value = "  +6.63%  "

$ws.Range("D31").Value = "'1.18"
$ws.Range("E31").Value = "  +0.88%  "

$ws.Range("D32").Value = "'3.22"
$ws.Range("E32").Value = "  +2.15%  "

$ws.Range("D33").Value = "'2.95"
$ws.Range("E33").Value = "  +1.00%  "

$ws.Range("E34").Value = "  +1.83%  "

$ws.Range("E35").Value = "  -0.58%  "

$ws.Range("D36").Value = "1.164.65"
$ws.Range("E36").Value = "  +0.94%  "

$ws.Range("D37").Value = "'0.0165"
$ws.Range("E37").Value = "  +0.20%  "

$ws.Range("E38").Value = "  +2.38%  "

$ws.Range("E39").Value = "  -0.22%  "

$ws.Range("E40").Value = "  -0.29%  "

$ws.Range("E41").Value = "  +1.12%  "

$ws.Range("D42").Value = "'5.43"
$ws.Range("E42").Value = "  +4.07%  "

$ws.Range("D43").Value = "'0.788"
$ws.Range("E43").Value = "  +1.30%  "

$ws.Range("D44").Value = "1.765.61"
$ws.Range("E44").Value = "  +1.58%  "

$ws.Range("D45").Value = "'92.67"
$ws.Range("E45").Value = "  +0.92%  "

$ws.Range("E46").Value = "  +3.00%  "

$ws.Range("D47").Value = "'54.48"
$ws.Range("E47").Value = "  +0.72%  "

$ws.Range("E48").Value = "  +1.08%  "

$ws.Range("D49").Value = "'0.409"
$ws.Range("E49").Value = "  +0.59%  "

$ws.Range("D50").Value = "'7.54"
$ws.Range("E50").Value = "  +5.02%  "

$ws.Range("E51").Value = "  -0.14%  "
